$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 35.906979
$ws.Range("H2").Value = 107.720937
$ws.Range("I2").Value = 0.6107087147789413
$ws.Range("J2").Value = 0.6107087147789412
$ws.Range("M2").Value = 6.977989333333333
$ws.Range("N2").Value = 20.933968
$ws.Range("O2").Value = 0.08060734116444308
$ws.Range("P2").Value = 0.0806073411644431
$ws.Range("Q2").Value = 250.558516454224
$ws.Range("R2").Value = 2255.026648088016
$ws.Range("S2").Value = 0.04922760572428468
$ws.Range("T2").Value = 0.04922760572428468

# Row 3
$ws.Range("G3").Value = 35.906979
$ws.Range("H3").Value = 107.720937
$ws.Range("I3").Value = 0.6107087147789413
$ws.Range("J3").Value = 0.6107087147789412
$ws.Range("O3").Value = 0.8031574845597853
$ws.Range("P3").Value = 0.8031574845597854
$ws.Range("Q3").Value = 2496.521345368165
$ws.Range("R3").Value = 22468.69210831348
$ws.Range("S3").Value = 0.4904952751605939
$ws.Range("T3").Value = 0.4904952751605938

# Row 4
$ws.Range("G4").Value = 35.906979
$ws.Range("H4").Value = 107.720937
$ws.Range("I4").Value = 0.6107087147789413
$ws.Range("J4").Value = 0.6107087147789412
$ws.Range("O4").Value = 0.1162351742757715
$ws.Range("P4").Value = 0.1162351742757716
$ws.Range("Q4").Value = 361.303479380639
$ws.Range("R4").Value = 3251.731314425751
$ws.Range("S4").Value = 0.0709858338940627
$ws.Range("T4").Value = 0.0709858338940627

# Row 5
$ws.Range("I5").Value = 0.2899643113254147
$ws.Range("J5").Value = 0.2899643113254147
$ws.Range("M5").Value = 6.977989333333333
$ws.Range("N5").Value = 20.933968
$ws.Range("O5").Value = 0.08060734116444308
$ws.Range("P5").Value = 0.0806073411644431
$ws.Range("Q5").Value = 118.9651071160249
$ws.Range("R5").Value = 1070.685964044224
$ws.Range("S5").Value = 0.02337325216852049
$ws.Range("T5").Value = 0.02337325216852049

# Row 6
$ws.Range("I6").Value = 0.2899643113254147
$ws.Range("J6").Value = 0.2899643113254147
$ws.Range("O6").Value = 0.8031574845597853
$ws.Range("P6").Value = 0.8031574845597854
$ws.Range("S6").Value = 0.2328870068962306
$ws.Range("T6").Value = 0.2328870068962305

# Row 7
$ws.Range("I7").Value = 0.2899643113254147
$ws.Range("J7").Value = 0.2899643113254147
$ws.Range("O7").Value = 0.1162351742757715
$ws.Range("P7").Value = 0.1162351742757716
$ws.Range("S7").Value = 0.03370405226066366
$ws.Range("T7").Value = 0.03370405226066366

# Row 8
$ws.Range("I8").Value = 0.09932697389564409
$ws.Range("J8").Value = 0.09932697389564407
$ws.Range("M8").Value = 6.977989333333333
$ws.Range("N8").Value = 20.933968
$ws.Range("O8").Value = 0.08060734116444308
$ws.Range("P8").Value = 0.0806073411644431
$ws.Range("Q8").Value = 40.75137397079467
$ws.Range("R8").Value = 366.762365737152
$ws.Range("S8").Value = 0.008006483271637915
$ws.Range("T8").Value = 0.008006483271637915

# Row 9
$ws.Range("I9").Value = 0.09932697389564409
$ws.Range("J9").Value = 0.09932697389564407
$ws.Range("O9").Value = 0.8031574845597853
$ws.Range("P9").Value = 0.8031574845597854
$ws.Range("S9").Value = 0.07977520250296097
$ws.Range("T9").Value = 0.07977520250296097

# Row 10
$ws.Range("I10").Value = 0.09932697389564409
$ws.Range("J10").Value = 0.09932697389564407
$ws.Range("O10").Value = 0.1162351742757715
$ws.Range("P10").Value = 0.1162351742757716
$ws.Range("Q10").Value = 58.76317202684135
$ws.Range("S10").Value = 0.0115452881210452
$ws.Range("T10").Value = 0.0115452881210452
